$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values in column D that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the sheet's original
# inline-string cell type) instead of silently converting them to numbers.

$ws.Range("D2").Value = "29.524.11"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.921.99"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'326.15"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("D7").Value = "'0.4808"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "'0.4042"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.08189"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "'1.007"
$ws.Range("D11").Value = "'23.79"
$ws.Range("D12").Value = "1.932.28"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'6.083"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'7.305"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "'91.52"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "'0.06874"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'1.011"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "29.526.43"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'11.98"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'2.177"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "2.166.50"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "'155.88"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").Value = "'6.386"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "'2.082"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "'120.40"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'1.013"
$ws.Range("D32").Value = "'0.09581"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'5.597"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").Value = "'3.561"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").Value = "'0.06351"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").Value = "'0.02284"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'1.191"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'0.5939"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'10.72"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'7.886"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'0.1841"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'2.480"
$ws.Range("E44").Value = "  +4.74%  "
$ws.Range("D45").Value = "'1.276"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'12.41"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").Value = "'0.5541"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "'1.972"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'117.76"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "'2.431"
$ws.Range("E51").Value = "  +0.94%  "

# Strip the forced-text number format so the cells stay style-index 0,
# matching the original (unstyled) cells; only their text content changed.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
